$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New trade row 7
$ws.Range("A7").Value = 9850.93
$ws.Range("B7").Value = 10217.75
$ws.Range("C7").Value = 296.89
$ws.Range("D7").Value = 286.23
$ws.Range("E7").Value = $false
$ws.Range("F7").Value = -3.59
$ws.Range("G7").Value = 42607.884317129632
$ws.Range("G7").NumberFormat = "m/d/yy h:mm"
$ws.Range("H7").Value = $false

# New trade row 8
$ws.Range("A8").Value = 9974.07
$ws.Range("B8").Value = 9850.93
$ws.Range("C8").Value = 286.39
$ws.Range("D8").Value = 282.82
$ws.Range("E8").Value = $true
$ws.Range("F8").Value = -1.25
$ws.Range("G8").Value = 42608.616388888891
$ws.Range("G8").NumberFormat = "m/d/yy h:mm"
$ws.Range("H8").Value = $true
